# edit.ps1 -- apply the spelling/grammar corrections described by the commit
# "Correction de quelque fautes d'orthographe dans le rapport"
#
# Implemented with Word's Find/Replace (wdReplaceAll behaviour via a single
# Execute call per fix) against $d.Content, which matches against the
# document's flattened text regardless of how it is split across runs.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute(
        $old,   # FindText
        $true,  # MatchCase
        $false, # MatchWholeWord
        $false, # MatchWildcards
        $false, # MatchSoundsLike
        $false, # MatchAllWordForms
        $true,  # Forward
        1,      # Wrap (wdFindContinue)
        $false, # Format
        $new,   # ReplaceWith
        2       # Replace (wdReplaceAll)
    ) | Out-Null
}

# 1. "Bombes qui détruisent ..." rule rewritten for clarity.
Replace-Text `
    "Bombes qui détruisent tous les bonbons avec lesquelles les bombes sont interverties." `
    "Bombes qui détruisent tous les bonbons de la même couleur que le bonbon avec lequel elles sont interverties."

# 2. "effectués" -> "effectué" (past participle agreement)
Replace-Text `
    "Après avoir effectués nos choix de langage et  " `
    "Après avoir effectué nos choix de langage et  "

# 3. "utilisez" -> "utilisé"
Replace-Text `
    "Pour finir, nous avons utilisez l’IDE " `
    "Pour finir, nous avons utilisé l’IDE "

# 4. "étés" -> "été"
Replace-Text `
    "sur Windows, ainsi tous les tests ont étés effectués sur Windows." `
    "sur Windows, ainsi tous les tests ont été effectués sur Windows."

# 5. "autre" -> "autres"
Replace-Text `
    ") et d’autre en privé" `
    ") et d’autres en privé"

# 6 & 7. "hérité" -> "héritée" (agreement with "classe Case") and
#        "case on créé" -> "case, on crée" (comma + verb conjugation)
Replace-Text `
    " est hérité de QQuickItem ce qui signifie que la vue et le modèle sont directement liés. En effet, pour créer une case on créé un composant qml " `
    " est héritée de QQuickItem ce qui signifie que la vue et le modèle sont directement liés. En effet, pour créer une case, on crée un composant qml "

# 8. "marqué" -> "marqués"
Replace-Text `
    "le score en fonction des bonbons marqué à détruire et " `
    "le score en fonction des bonbons marqués à détruire et "

# 9. "possible" -> "possibles"
Replace-Text `
    "si il reste des coups possible et appelle " `
    "si il reste des coups possibles et appelle "

# 10. "permit" -> "permis"
Replace-Text `
    "ce fut une découverte. Cela nous également permit d’apprendre" `
    "ce fut une découverte. Cela nous également permis d’apprendre"
